# Auto-push with fireball accuracy report [2025-04-15 04:18 AM]
# Update the ATS Accuracy and Total Accuracy sheets with refreshed
# fireball accuracy stats.

$wb = $excel.ActiveWorkbook

# --- Sheet: ATS Accuracy ---
$wsATS = $wb.Worksheets.Item("ATS Accuracy")

# Row 2 (5 fireballs)
$wsATS.Range("C2").Value = 83
$wsATS.Range("D2").Value = 86
$wsATS.Range("E2").Value = 96.5

# Row 3 (4 fireballs)
$wsATS.Range("B3").Value = 3
$wsATS.Range("D3").Value = 62
$wsATS.Range("E3").Value = 95.2

# Row 4 (3 fireballs)
$wsATS.Range("C4").Value = 14
$wsATS.Range("D4").Value = 17
$wsATS.Range("E4").Value = 82.40000000000001

# Row 5 (2 fireballs)
$wsATS.Range("B5").Value = 2
$wsATS.Range("C5").Value = 7
$wsATS.Range("D5").Value = 9
$wsATS.Range("E5").Value = 77.8

# Row 6 (1 fireball)
$wsATS.Range("C6").Value = 4
$wsATS.Range("D6").Value = 8
$wsATS.Range("E6").Value = 50

# --- Sheet: Total Accuracy ---
$wsTotal = $wb.Worksheets.Item("Total Accuracy")

# Row 2 (5 fireballs)
$wsTotal.Range("B2").Value = 5
$wsTotal.Range("C2").Value = 72
$wsTotal.Range("D2").Value = 77
$wsTotal.Range("E2").Value = 93.5

# Row 3 (4 fireballs)
$wsTotal.Range("B3").Value = 5
$wsTotal.Range("C3").Value = 61
$wsTotal.Range("E3").Value = 92.40000000000001

# Row 4 (3 fireballs)
$wsTotal.Range("B4").Value = 1
$wsTotal.Range("C4").Value = 23
$wsTotal.Range("D4").Value = 24
$wsTotal.Range("E4").Value = 95.8

# Row 5 (2 fireballs) - unchanged

# Row 6 (1 fireball)
$wsTotal.Range("B6").Value = 4
$wsTotal.Range("C6").Value = 0
$wsTotal.Range("D6").Value = 4
$wsTotal.Range("E6").Value = 0
